$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the test case result columns (G = Population Change, H = Future Population)
$ws.Range("G6").Value = 7821428
$ws.Range("H6").Value = 340921788

$ws.Range("G7").Value = 19647428
$ws.Range("H7").Value = 352747788

# Row 8: correct the "No of Years" value in D8, and fill G8/H8
$ws.Range("D8").Value = 80
$ws.Range("G8").Value = 8541000
$ws.Range("H8").Value = 341641360

# Row 9: only Future Population filled in
$ws.Range("H9").Value = 334351788

# Row 10: only Future Population filled in
$ws.Range("H10").Value = 332011617

# Update the active selection on the sheet
$ws.Range("G10").Select()
